$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
# Swap Linea values between rows 35 and 36
$ws.Cells.Item(35,3).Value = "14_ABASTO"
$ws.Cells.Item(36,3).Value = "10_OLMOS"
# Update header text
$ws.Cells.Item(2,1).Value = "Última actualización: 08:50:00"
$ws.Cells.Item(3,1).Value = "Total filas: 126"
# Insert new scraped rows (ascending order so indices remain valid)
$ws.Rows.Item(87).Insert()
$ws.Cells.Item(87,1).Value = "08:50:00"
$ws.Cells.Item(87,2).Value = "08:50"
$ws.Cells.Item(87,3).Value = "10_OLMOS"
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = "LP1912"
$ws.Rows.Item(101).Insert()
$ws.Cells.Item(101,1).Value = "08:50:00"
$ws.Cells.Item(101,2).Value = "09:13"
$ws.Cells.Item(101,3).Value = "10_OLMOS"
$ws.Cells.Item(101,4).Value = 23
$ws.Cells.Item(101,5).Value = "LP1912"
$ws.Rows.Item(114).Insert()
$ws.Cells.Item(114,1).Value = "08:50:00"
$ws.Cells.Item(114,2).Value = "09:33"
$ws.Cells.Item(114,3).Value = "16_SANTA ANA"
$ws.Cells.Item(114,4).Value = 43
$ws.Cells.Item(114,5).Value = "LP1912"
$ws.Rows.Item(118).Insert()
$ws.Cells.Item(118,1).Value = "08:50:00"
$ws.Cells.Item(118,2).Value = "09:35"
$ws.Cells.Item(118,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(118,4).Value = 45
$ws.Cells.Item(118,5).Value = "LP1912"
$ws.Rows.Item(120).Insert()
$ws.Cells.Item(120,1).Value = "08:50:00"
$ws.Cells.Item(120,2).Value = "09:38"
$ws.Cells.Item(120,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(120,4).Value = 48
$ws.Cells.Item(120,5).Value = "LP1912"
$ws.Rows.Item(129).Insert()
$ws.Cells.Item(129,1).Value = "08:50:00"
$ws.Cells.Item(129,2).Value = "10:27"
$ws.Cells.Item(129,3).Value = "215A_EL PATO"
$ws.Cells.Item(129,4).Value = 97
$ws.Cells.Item(129,5).Value = "LP1912"
$ws.Rows.Item(130).Insert()
$ws.Cells.Item(130,1).Value = "08:50:00"
$ws.Cells.Item(130,2).Value = "10:42"
$ws.Cells.Item(130,3).Value = "17_ROMERO"
$ws.Cells.Item(130,4).Value = 112
$ws.Cells.Item(130,5).Value = "LP1912"
$ws.Rows.Item(131).Insert()
$ws.Cells.Item(131,1).Value = "08:50:00"
$ws.Cells.Item(131,2).Value = "10:44"
$ws.Cells.Item(131,3).Value = "14_ABASTO"
$ws.Cells.Item(131,4).Value = 114
$ws.Cells.Item(131,5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
# Update header text
$ws.Cells.Item(2,1).Value = "Última actualización: 08:50:00"
$ws.Cells.Item(3,1).Value = "Total filas: 19"
# Insert new scraped rows (ascending order so indices remain valid)
$ws.Rows.Item(24).Insert()
$ws.Cells.Item(24,1).Value = "08:50:00"
$ws.Cells.Item(24,2).Value = "10:27"
$ws.Cells.Item(24,3).Value = "215A_EL PATO"
$ws.Cells.Item(24,4).Value = 97
$ws.Cells.Item(24,5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
# Update header text
$ws.Cells.Item(2,1).Value = "Última actualización: 08:50:00"
$ws.Cells.Item(3,1).Value = "Total filas: 27"
# Insert new scraped rows (ascending order so indices remain valid)
$ws.Rows.Item(28).Insert()
$ws.Cells.Item(28,1).Value = "08:50:00"
$ws.Cells.Item(28,2).Value = "08:58"
$ws.Cells.Item(28,3).Value = "215A_LA PLATA"
$ws.Cells.Item(28,4).Value = 8
$ws.Cells.Item(28,5).Value = "L6173"
